# fix: unique command names in XLSX - prefix protocol name to each step
#
# For every worksheet whose name matches one of the "protocol" sheets
# (i.e. every sheet except the five overview/profile sheets at the
# front of the workbook), prefix each Step/command name in column A
# (rows 2..last used row) with "<SheetName> " so that command names are
# globally unique across sheets.

$wb = $excel.ActiveWorkbook

# These leading sheets hold general profile/journey data, not per-sheet
# "commands" - they are left untouched.
$skip = @("VeraJourney", "NRWaves", "PersonalVera", "PositiveSpin", "ReEngagement")

foreach ($ws in $wb.Worksheets) {
    if ($skip -contains $ws.Name) {
        continue
    }

    $lastRow = $ws.UsedRange.Rows.Count

    for ($r = 2; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, 1)
        $current = $cell.Value()
        if ($current -ne $null -and $current -ne "") {
            $prefix = $ws.Name + " "
            if (-not $current.StartsWith($prefix)) {
                $cell.Value = $prefix + $current
            }
        }
    }
}
